# Applies the cell-level corrections described in the commit diff to the
# "solar" and "wind" sheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- solar sheet: swap lcoe_class (P) values between rows 30 and 31 ---
$wsSolar = $wb.Worksheets.Item("solar")
$wsSolar.Range("P30").Value = 3
$wsSolar.Range("P31").Value = 4

# --- wind sheet: correct cap_bnd (M), ncap_cost (O) and lcoe_class (P) ---
$wsWind = $wb.Worksheets.Item("wind")

# Rows 7 and 8: swap cap_bnd, ncap_cost and lcoe_class values
$wsWind.Range("M7").Value = 0.0015
$wsWind.Range("O7").Value = 31.372355844942916
$wsWind.Range("P7").Value = 1

$wsWind.Range("M8").Value = 2.4990000000000001
$wsWind.Range("O8").Value = 39.630069093581724
$wsWind.Range("P8").Value = 3

# Rows 14, 15, 16: rotate lcoe_class values
$wsWind.Range("P14").Value = 4
$wsWind.Range("P15").Value = 5
$wsWind.Range("P16").Value = 3

# Rows 24, 25: swap lcoe_class values
$wsWind.Range("P24").Value = 1
$wsWind.Range("P25").Value = 2

# Rows 98, 99: swap lcoe_class values
$wsWind.Range("P98").Value = 1
$wsWind.Range("P99").Value = 2
